{"js": "// Insert the word \"i\" (with surrounding space) into the first paragraph,\n// turning:\n//   Svar p\u00e5 \u00abLeveranse: Innhold i dokumentet\u00bb - eksamensteksten.\n// into:\n//   Svar p\u00e5 \u00abLeveranse: Innhold i dokumentet\u00bb - i eksamensteksten.\n\nconst body = context.document.body;\n\n// Locate the unique anchor text right after the dash, where the new\n// word \"i\" needs to be inserted.\nconst results = body.search(\" eksamensteksten.\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Anchor text \" eksamensteksten.\" not found.');\n}\n\n// Insert \" i\" immediately before \" eksamensteksten.\" (i.e. right after the\n// trailing \"-\" of the previous sentence), producing \"- i eksamensteksten.\".\nconst anchor = results.items[0];\nanchor.insertText(\" i\", Word.InsertLocation.before);\n\nawait context.sync();\n", "ps1": "# Insert the word \"i\" (with surrounding space) into the first paragraph,\n# turning:\n#   Svar p\u00e5 \u00abLeveranse: Innhold i dokumentet\u00bb - eksamensteksten.\n# into:\n#   Svar p\u00e5 \u00abLeveranse: Innhold i dokumentet\u00bb - i eksamensteksten.\n\n$d = $word.ActiveDocument\n\n# Find the unique anchor text right after the dash, where the new word\n# \"i\" needs to be inserted.\n$rng = $d.Content\n$find = $rng.Find\n$find.Text = \" eksamensteksten.\"\n$found = $find.Execute()\n\nif (-not $found) {\n    throw 'Anchor text \" eksamensteksten.\" not found.'\n}\n\n# $rng now covers \" eksamensteksten.\"; collapse to its start and insert\n# \" i\" right before it, producing \"- i eksamensteksten.\".\n$rng.Collapse(1)  # wdCollapseStart\n$rng.InsertBefore(\" i\")\n"}
